$wb = $excel.ActiveWorkbook

# ----- Sheet1: Neg_Change -----
$ws1 = $wb.Worksheets.Item("Neg_Change")

$arr1 = New-Object "object[,]" 15,8
$arr1[0,0] = "LT"
$arr1[0,1] = 4072.8
$arr1[0,2] = 4096
$arr1[0,3] = 4063.8
$arr1[0,4] = 4090
$arr1[0,5] = 1064975
$arr1[0,6] = 2207898
$arr1[0,7] = -0.5176520835654546
$arr1[1,0] = "ULTRACEMCO"
$arr1[1,1] = 11658
$arr1[1,2] = 11804
$arr1[1,3] = 11658
$arr1[1,4] = 11711
$arr1[1,5] = 178073
$arr1[1,6] = 389481
$arr1[1,7] = -0.5427941285967737
$arr1[2,0] = "AMBUJACEM"
$arr1[2,1] = 547.75
$arr1[2,2] = 556.3
$arr1[2,3] = 542.95
$arr1[2,4] = 554
$arr1[2,5] = 1239318
$arr1[2,6] = 2488693
$arr1[2,7] = -0.5020205384914893
$arr1[3,0] = "ADANIENSOL"
$arr1[3,1] = 1014.4
$arr1[3,2] = 1020
$arr1[3,3] = 1007.05
$arr1[3,4] = 1015.2
$arr1[3,5] = 887627
$arr1[3,6] = 2110784
$arr1[3,7] = -0.5794799467875443
$arr1[4,0] = "JINDALSTEL"
$arr1[4,1] = 1032
$arr1[4,2] = 1036.5
$arr1[4,3] = 1023.2
$arr1[4,4] = 1031
$arr1[4,5] = 396392
$arr1[4,6] = 846932
$arr1[4,7] = -0.5319671473034435
$arr1[5,0] = "TATAPOWER"
$arr1[5,1] = 380
$arr1[5,2] = 382.4
$arr1[5,3] = 378
$arr1[5,4] = 382
$arr1[5,5] = 3359948
$arr1[5,6] = 6778054
$arr1[5,7] = -0.5042901694203086
$arr1[6,0] = "BOSCHLTD"
$arr1[6,1] = 36870
$arr1[6,2] = 36870
$arr1[6,3] = 36340
$arr1[6,4] = 36390
$arr1[6,5] = 4956
$arr1[6,6] = 11065
$arr1[6,7] = -0.5521012200632626
$arr1[7,0] = "UPL"
$arr1[7,1] = 751
$arr1[7,2] = 769.55
$arr1[7,3] = 743
$arr1[7,4] = 766
$arr1[7,5] = 1840513
$arr1[7,6] = 4075773
$arr1[7,7] = -0.5484260286330961
$arr1[8,0] = "SONACOMS"
$arr1[8,1] = 490
$arr1[8,2] = 496
$arr1[8,3] = 483.1
$arr1[8,4] = 485.5
$arr1[8,5] = 1169034
$arr1[8,6] = 2896193
$arr1[8,7] = -0.5963549390527496
$arr1[9,0] = "HDFCAMC"
$arr1[9,1] = 2669.8
$arr1[9,2] = 2669.8
$arr1[9,3] = 2597.2
$arr1[9,4] = 2608
$arr1[9,5] = 626759
$arr1[9,6] = 1304915
$arr1[9,7] = -0.5196936198909508
$arr1[10,0] = "KFINTECH"
$arr1[10,1] = 1056
$arr1[10,2] = 1056.9
$arr1[10,3] = 1039.2
$arr1[10,4] = 1049.7
$arr1[10,5] = 402076
$arr1[10,6] = 837476
$arr1[10,7] = -0.5198954955127072
$arr1[11,0] = "GRANULES"
$arr1[11,1] = 576.95
$arr1[11,2] = 577.75
$arr1[11,3] = 569.4
$arr1[11,4] = 571.9
$arr1[11,5] = 583309
$arr1[11,6] = 1274628
$arr1[11,7] = -0.5423692245894488
$arr1[12,0] = "PGEL"
$arr1[12,1] = 566
$arr1[12,2] = 577
$arr1[12,3] = 561.3
$arr1[12,4] = 564
$arr1[12,5] = 1825949
$arr1[12,6] = 4093056
$arr1[12,7] = -0.55389102909904
$arr1[13,0] = "MANAPPURAM"
$arr1[13,1] = 287.55
$arr1[13,2] = 290.25
$arr1[13,3] = 285.1
$arr1[13,4] = 285.3
$arr1[13,5] = 1150976
$arr1[13,6] = 2857546
$arr1[13,7] = -0.5972152329306335
$arr1[14,0] = "NUVAMA"
$arr1[14,1] = 7349
$arr1[14,2] = 7349
$arr1[14,3] = 7242.5
$arr1[14,4] = 7275
$arr1[14,5] = 33266
$arr1[14,6] = 72775
$arr1[14,7] = -0.5428924768120921
$ws1.Range("A2:H16").Value = $arr1

# Mirror symbol into column I (filterdata)
$arrI1 = New-Object "object[,]" 15,1
$arrI1[0,0] = "LT"
$arrI1[1,0] = "ULTRACEMCO"
$arrI1[2,0] = "AMBUJACEM"
$arrI1[3,0] = "ADANIENSOL"
$arrI1[4,0] = "JINDALSTEL"
$arrI1[5,0] = "TATAPOWER"
$arrI1[6,0] = "BOSCHLTD"
$arrI1[7,0] = "UPL"
$arrI1[8,0] = "SONACOMS"
$arrI1[9,0] = "HDFCAMC"
$arrI1[10,0] = "KFINTECH"
$arrI1[11,0] = "GRANULES"
$arrI1[12,0] = "PGEL"
$arrI1[13,0] = "MANAPPURAM"
$arrI1[14,0] = "NUVAMA"
$ws1.Range("I2:I16").Value = $arrI1

# ----- Sheet2: Pos_Change -----
$ws2 = $wb.Worksheets.Item("Pos_Change")

$arr2 = New-Object "object[,]" 9,8
$arr2[0,0] = "BAJAJ-AUTO"
$arr2[0,1] = 9040
$arr2[0,2] = 9042.5
$arr2[0,3] = 8884
$arr2[0,4] = 8921
$arr2[0,5] = 507849
$arr2[0,6] = 356060
$arr2[0,7] = 0.4263017468965904
$arr2[1,0] = "M&M"
$arr2[1,1] = 3632
$arr2[1,2] = 3665.4
$arr2[1,3] = 3603.6
$arr2[1,4] = 3609.7
$arr2[1,5] = 1704069
$arr2[1,6] = 1078982
$arr2[1,7] = 0.5793303317386203
$arr2[2,0] = "PIDILITIND"
$arr2[2,1] = 1480
$arr2[2,2] = 1483.4
$arr2[2,3] = 1468.4
$arr2[2,4] = 1475.3
$arr2[2,5] = 494291
$arr2[2,6] = 332259
$arr2[2,7] = 0.4876677531684619
$arr2[3,0] = "LICI"
$arr2[3,1] = 864
$arr2[3,2] = 866.15
$arr2[3,3] = 855.5
$arr2[3,4] = 856.55
$arr2[3,5] = 760013
$arr2[3,6] = 513683
$arr2[3,7] = 0.4795369907121708
$arr2[4,0] = "NAUKRI"
$arr2[4,1] = 1362.6
$arr2[4,2] = 1371.7
$arr2[4,3] = 1345.3
$arr2[4,4] = 1360
$arr2[4,5] = 932411
$arr2[4,6] = 594311
$arr2[4,7] = 0.5688940638823781
$arr2[5,0] = "COLPAL"
$arr2[5,1] = 2161
$arr2[5,2] = 2171.5
$arr2[5,3] = 2146.1
$arr2[5,4] = 2163.4
$arr2[5,5] = 396185
$arr2[5,6] = 248587
$arr2[5,7] = 0.5937478629212308
$arr2[6,0] = "MPHASIS"
$arr2[6,1] = 2895.7
$arr2[6,2] = 2903.5
$arr2[6,3] = 2860.3
$arr2[6,4] = 2895
$arr2[6,5] = 207272
$arr2[6,6] = 140084
$arr2[6,7] = 0.4796265098084007
$arr2[7,0] = "INDIANB"
$arr2[7,1] = 788
$arr2[7,2] = 793.25
$arr2[7,3] = 779.05
$arr2[7,4] = 783
$arr2[7,5] = 2416248
$arr2[7,6] = 1533875
$arr2[7,7] = 0.575257436231766
$arr2[8,0] = "CESC"
$arr2[8,1] = 169
$arr2[8,2] = 171.5
$arr2[8,3] = 168.64
$arr2[8,4] = 171.09
$arr2[8,5] = 604975
$arr2[8,6] = 428230
$arr2[8,7] = 0.4127338112696448
$ws2.Range("A2:H10").Value = $arr2

# Mirror symbol into column I (filterdata)
$arrI2 = New-Object "object[,]" 9,1
$arrI2[0,0] = "BAJAJ-AUTO"
$arrI2[1,0] = "M&M"
$arrI2[2,0] = "PIDILITIND"
$arrI2[3,0] = "LICI"
$arrI2[4,0] = "NAUKRI"
$arrI2[5,0] = "COLPAL"
$arrI2[6,0] = "MPHASIS"
$arrI2[7,0] = "INDIANB"
$arrI2[8,0] = "CESC"
$ws2.Range("I2:I10").Value = $arrI2

# Sheet2 shrank from 13 data rows (12 entries) to 10 rows (9 entries); delete leftover rows
$ws2.Range("A11:I13").EntireRow.Delete()

